# Commit: "Added Display class Updated images"
#
# The "IsNew" and "IsTopBuyed" boolean columns (G, H) are removed and
# replaced by a single new "DisplayClass" column in their place (G).
# All following columns shift left by one (U -> T dimension).
# Selection moves from D5:D7 to G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "IsTopBuyed" (G) and "IsNew" (H) columns entirely.
$ws.Range("G1:H1").EntireColumn.Delete()

# Insert a new empty column in their place for "DisplayClass".
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").Value = "DisplayClass"

# Match the author's final selection (cell G2).
$ws.Range("G2").Select() | Out-Null
